# Weekly data update: a new weekly record is inserted as row 78, pushing
# all subsequent rows (old 78-188) down by one (new 79-189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 - shifts rows 78..188 down to 79..189
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with this week's record
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 44546
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100103
$ws.Range("H78").Value = "Frutos de hueso (carozo)"
$ws.Range("I78").Value = 100103004
$ws.Range("J78").Value = "Durazno"
$ws.Range("K78").Value = "Royal Glory"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 120
$ws.Range("N78").Value = 17000
$ws.Range("O78").Value = 18000
$ws.Range("P78").Value = 17500
$ws.Range("Q78").Value = "$/caja 16 kilos empedrada"
$ws.Range("R78").Value = "Región de O'Higgins"
$ws.Range("S78").Value = 1094
$ws.Range("T78").Value = 16
